$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the contact values (FirstName / Street / Hobbies) ---
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"

$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Normalize the pincode/phone column font from the theme color to explicit black ---
$ws.Range("I2").Font.Color = 0
$ws.Range("K2").Font.Color = 0
$ws.Range("I3").Font.Color = 0
$ws.Range("K3").Font.Color = 0

# --- Row heights bumped slightly on re-save ---
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
